{"js": "// Highlight the key criteria names in green for the two checklist items:\n//   \"Horas continuas/dia de un profesor\"  -> highlight \"Horas continuas\"\n//   \"M\u00e1x de Gaps por d\u00eda\"                 -> highlight \"M\u00e1x de Gaps\"\n// The trailing \"/\" and \" por d\u00eda\" stay un-highlighted (bold only).\n// Moving the highlight boundary on the \"M\u00e1x de Gaps\" item also relocates\n// the document's `_GoBack` bookmark (Word's \"last edit\" marker) to sit\n// right after that edit, matching the author's original Word session.\n\nconst body = context.document.body;\n\n// --- 1) \"Horas continuas/dia de un profesor\" ---------------------------\nconst continuas = body.search(\"Horas continuas\", { matchCase: true });\ncontinuas.load(\"text\");\nawait context.sync();\n\nif (continuas.items.length > 0) {\n  continuas.items[0].font.highlightColor = \"green\";\n}\n\n// --- 2) \"M\u00e1x de Gaps por d\u00eda\" (not the \"...por semana\" sibling) --------\nconst gapsPorDia = body.search(\"M\u00e1x de Gaps por d\u00eda\", { matchCase: true });\ngapsPorDia.load(\"text\");\nawait context.sync();\n\nif (gapsPorDia.items.length > 0) {\n  const fullItem = gapsPorDia.items[0];\n\n  // Highlight just the \"M\u00e1x de Gaps\" portion of that match.\n  const target = fullItem.search(\"M\u00e1x de Gaps\", { matchCase: true });\n  target.load(\"text\");\n  await context.sync();\n\n  if (target.items.length > 0) {\n    target.items[0].font.highlightColor = \"green\";\n\n    // Relocate the `_GoBack` bookmark to the end of the newly-highlighted\n    // text (right before \" por d\u00eda\"), removing it from its old location.\n    const insertionPoint = target.items[0].getRange(\"End\");\n    await context.sync();\n\n    context.document.deleteBookmark(\"_GoBack\");\n    insertionPoint.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Highlight the key criteria names in green for the two checklist items:\n#   \"Horas continuas/dia de un profesor\"  -> highlight \"Horas continuas\"\n#   \"M\u00e1x de Gaps por d\u00eda\"                 -> highlight \"M\u00e1x de Gaps\"\n# The trailing \"/\" and \" por d\u00eda\" stay un-highlighted (bold only).\n# Moving the highlight boundary on the \"M\u00e1x de Gaps\" item also relocates\n# the document's `_GoBack` bookmark (Word's \"last edit\" marker) to sit\n# right after that edit, matching the author's original Word session.\n#\n# Note: on this host, Range.HighlightColorIndex (set directly on a Range)\n# does not reliably target the Range; Range.Font.HighlightColorIndex does,\n# so that is what is used below. 4 = wdBrightGreen, which maps to the\n# OOXML <w:highlight w:val=\"green\"/> value used by the target edit.\n\n$d = $word.ActiveDocument\n\n# --- 1) \"Horas continuas/dia de un profesor\" ----------------------------\n$r1 = $d.Content\n$r1.Find.MatchCase = $true\n$found1 = $r1.Find.Execute(\"Horas continuas\")\nif ($found1) {\n    $r1.Font.HighlightColorIndex = 4\n}\n\n# --- 2) \"M\u00e1x de Gaps por d\u00eda\" (not the \"...por semana\" sibling) --------\n$r2 = $d.Content\n$r2.Find.MatchCase = $true\n$found2 = $r2.Find.Execute(\"M\u00e1x de Gaps por d\u00eda\")\nif ($found2) {\n    $highlightText = \"M\u00e1x de Gaps\"\n    $target = $d.Range($r2.Start, $r2.Start + $highlightText.Length)\n    $target.Font.HighlightColorIndex = 4\n\n    # Relocate the `_GoBack` bookmark to the end of the newly-highlighted\n    # text (right before \" por d\u00eda\"). Adding a bookmark with an existing\n    # name moves it, so any prior `_GoBack` elsewhere in the document is\n    # removed automatically.\n    $insertionPoint = $d.Range($target.End, $target.End)\n    $d.Bookmarks.Add(\"_GoBack\", $insertionPoint)\n}\n"}
